$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" row (row 26) entirely.
$ws.Rows(26).Delete()

# After the first deletion, "SC 92" (originally row 28) is now row 27. Remove it too.
$ws.Rows(27).Delete()

# Fill in previously-missing value for "SC 5" (now row 26), column D.
$ws.Range("D26").Value = -13.8

# Clear value for "SC 101" (now row 27), column D (now treated as missing).
$ws.Range("D27").ClearContents()

# Fill in previously-missing value for "SC 232" (now row 33), column D.
$ws.Range("D33").Value = -14.1
